$d = $word.ActiveDocument
$p3 = $d.Paragraphs(3)

# Step 0: set entire paragraph (mark + text) to FFFF00 first (this colors pPr mark and the single run)
$p3.Range.Font.Color = 65535

$pStart = $p3.Range.Start
$pEndNoMark = $p3.Range.End - 1
$full = $d.Range($pStart, $pEndNoMark).Text

# Step 1: Replace "ceptable" with "preciable"
$idxCeptable = $full.IndexOf("ceptable")
$ceptStart = $pStart + $idxCeptable
$ceptEnd = $ceptStart + 8
$d.Range($ceptStart, $ceptEnd).Text = "preciable"

$preciableStart = $ceptStart
$preciableEnd = $ceptStart + 9
$periodStart = $preciableEnd
$periodEnd = $periodStart + 1

# Step 2: remove the _GoBack bookmark
$d.Bookmarks("_GoBack").Delete()

# Step 3: insert new sentence after the period (inherits current FFFF00, will fix colors below)
$insertPoint = $d.Range($periodEnd, $periodEnd)
$insertPoint.InsertAfter(" Sin embargo estos dispositivos no responden a sobretemperaturas causadas por condiciones ambientales (temperatura ambiente superior a 40°C y por fallas de ventilación).")

$spaceStart = $periodEnd
$spaceEnd = $spaceStart + 1
$sentStart = $spaceEnd
$sentEnd = $p3.Range.End - 1

# Step 4: carve out correct colors by toggling to a different value first, to force run splits
$rLead = $d.Range($pStart, $preciableStart)
$rLead.Font.Color = 65535
$rLead.Font.Color = 255

$rPreciable = $d.Range($preciableStart, $preciableEnd)
$rPreciable.Font.Color = 65535
$rPreciable.Font.Color = 255

$rPeriod = $d.Range($periodStart, $periodEnd)
$rPeriod.Font.Color = 65535
$rPeriod.Font.Color = 255

$rSpace = $d.Range($spaceStart, $spaceEnd)
$rSpace.Font.Color = 65535
$rSpace.Font.Color = 255

$rSent = $d.Range($sentStart, $sentEnd)
$rSent.Font.Color = 255
$rSent.Font.Color = 65535

# Step 5: insert a new empty paragraph after paragraph 3, with pPr/rPr color FFFF00 + lang es-ES,
# using InsertXML to avoid an unwanted empty run artifact.
$endOfP3 = $p3.Range.End
$insertParaPoint = $d.Range($endOfP3, $endOfP3)
$xmlFragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:color w:val="FFFF00"/><w:lang w:val="es-ES"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertParaPoint.InsertXML($xmlFragment)

Write-Output "DONE"
